# Fruta / hortaliza, semanal
# Insert the latest weekly price observation at the top of the data table
# (row 169), pushing all existing data rows down by one and extending the
# used range from A1:R228 to A1:R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right before the current first data row (row 169),
# shifting rows 169:228 down to 170:229.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A169").Value = 5
$ws.Range("B169").Value = "Macroferia Regional de Talca"
$ws.Range("C169").Value = "Maule"
$ws.Range("D169").Value = 44985
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 100112031
$ws.Range("G169").Value = "Poroto verde"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 150
$ws.Range("K169").Value = 25000
$ws.Range("L169").Value = 25000
$ws.Range("M169").Value = 25000
$ws.Range("N169").Value = "`$/saco 25 kilos"
$ws.Range("O169").Value = "Región del Maule"
$ws.Range("P169").Value = 1000
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"
